# Generate Report for Handoff
# Refresh the handoff timestamps for the "b04048be-ba5b-487b-a279-7de58c31591f.md"
# file (row 5 of the per-locale sheets) across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2017-02-09 07:53:47"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2017-02-09 07:53:29"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2017-02-09 07:53:47"
